# risk_factors.xlsx - update factor_description wording (column C) for every
# factor row on the risk_factors sheet, then refresh the column widths /
# selection to reflect the edited content, mirroring the authored commit
# ("Modified product types labels in script and and factors questions labels
# in risk factors excel file").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- inherent risk factors ---------------------------------------------
$ws.Range("C2").Value = "Direct distribution to vulnerable populations (e.g., hospitals, elderly care facilities, neonatal units)"
$ws.Range("C3").Value = "Food specifically intended for vulnerable populations (e.g., infant formula, meals for elderly individuals)"

# --- mitigation factors ---------------------------------------------------
$ws.Range("C4").Value  = "Certification or accreditation of the business food safety management system (e.g. HACCP, ISO 22000) by a reputable body"
$ws.Range("C5").Value  = "Adequate food safety training for employees (e.g., PAHO Food Handlers Manual), especially those handling or processing food"
$ws.Range("C6").Value  = "Official recognition by competent authorities for export eligibility"
$ws.Range("C7").Value  = "Verified access to potable water (e.g., annual on-site tests or municipal certification)"
$ws.Range("C8").Value  = "Functional and well-maintained sanitary facilities available for personnel (e.g., toilets)"
$ws.Range("C9").Value  = "Functional and well-maintained handwashing facilities accessible to personnel"
$ws.Range("C10").Value = "Reliable and continuous access to electrical power"
$ws.Range("C11").Value = "Adequate and continuously available cooling equipment based on product requirements (e.g., refrigeration for perishable goods)"

# --- compliance factors ---------------------------------------------------
$ws.Range("C12").Value = "Current inspection identified at least one major food safety non-conformity"
$ws.Range("C13").Value = "The food business has faced enforcement actions by inspection authorities within the past 5 years (e.g., permit suspension, temporary closure)"
$ws.Range("C14").Value = "Previous inspection identified at least one major food safety non-conformity"
$ws.Range("C15").Value = "Food business linked to a documented foodborne outbreak within the past 3 years"

# --- resize the columns so the (now differently-sized) text is fully
# visible, closest achievable widths to A=10.44, B=19.44, C=155.11 chars ---
$ws.Columns.Item(1).ColumnWidth = 9.666666666666666
$ws.Columns.Item(2).ColumnWidth = 18.666666666666668
$ws.Columns.Item(3).ColumnWidth = 154.33333333333331

# --- leave the selection on C19, matching the saved view state -----------
$ws.Range("C19").Select() | Out-Null

Write-Output "risk_factors descriptions updated"
